# Se eliminan los Estados de Cuenta (EC) anteriores y se agregan nuevos;
# se modifica la base de datos de la hoja "Hoja1".
#
# La tabla de deudas (filas 16-20, columnas B:G) se reorganiza: las filas
# correspondientes a los periodos de mora de ISIDORO SALAS AGUIRRE quedan en
# orden descendente (1807, 1806, 1805, 1804) y el registro de PATRICIA
# MARTINEZ MUTIS (periodo 1804) pasa a ser el ultimo de la tabla.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 16: ISIDORO SALAS AGUIRRE - periodo 1807
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "15609536"
$ws.Range("D16").Value = "ISIDORO SALAS AGUIRRE"
$ws.Range("E16").Value = "1807"
$ws.Range("F16").Value = 72000
$ws.Range("G16").Value = 1800000

# Fila 17: ISIDORO SALAS AGUIRRE - periodo 1806
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "15609536"
$ws.Range("D17").Value = "ISIDORO SALAS AGUIRRE"
$ws.Range("E17").Value = "1806"
$ws.Range("F17").Value = 72000
$ws.Range("G17").Value = 1800000

# Fila 18: ISIDORO SALAS AGUIRRE - periodo 1805 (sin cambios)
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "15609536"
$ws.Range("D18").Value = "ISIDORO SALAS AGUIRRE"
$ws.Range("E18").Value = "1805"
$ws.Range("F18").Value = 72000
$ws.Range("G18").Value = 1800000

# Fila 19: ISIDORO SALAS AGUIRRE - periodo 1804
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "15609536"
$ws.Range("D19").Value = "ISIDORO SALAS AGUIRRE"
$ws.Range("E19").Value = "1804"
$ws.Range("F19").Value = 69600
$ws.Range("G19").Value = 1800000

# Fila 20: PATRICIA MARTINEZ MUTIS - periodo 1804
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "45488904"
$ws.Range("D20").Value = "PATRICIA MARTINEZ MUTIS"
$ws.Range("E20").Value = "1804"
$ws.Range("F20").Value = 17709
$ws.Range("G20").Value = 781300
